# Auto-generated edit script applying numeric updates to the
# "Excalibur_Profits" market-price/profit columns (H:N) across
# the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (46 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 430.9524
$ws.Range("I28").Value = 421.8421
$ws.Range("J28").Value = 517.5
$ws.Range("K28").Value = 421.8421
$ws.Range("L28").Value = 517.5
$ws.Range("M28").Value = 63.15789999999998
$ws.Range("N28").Value = -1487.5
$ws.Range("H51").Value = 2000.5
$ws.Range("I51").Value = 2000.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2000.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1516.5
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 513.6
$ws.Range("J58").Value = 401
$ws.Range("L58").Value = 1203
$ws.Range("N58").Value = -1503
$ws.Range("H61").Value = 967.4286
$ws.Range("I61").Value = 967.4286
$ws.Range("K61").Value = 2902.2858
$ws.Range("M61").Value = -2730.2858
$ws.Range("H106").Value = 3819.5454
$ws.Range("I106").Value = 3832.4285
$ws.Range("K106").Value = 3832.4285
$ws.Range("M106").Value = -3201.4285
$ws.Range("H112").Value = 10598.296
$ws.Range("J112").Value = 12843.546
$ws.Range("L112").Value = 38530.638
$ws.Range("N112").Value = -40746.638
$ws.Range("H127").Value = 3597.1428
$ws.Range("I127").Value = 3597.1428
$ws.Range("K127").Value = 10791.4284
$ws.Range("M127").Value = -5831.428400000001
$ws.Range("H132").Value = 77640.836
$ws.Range("I132").Value = 88321.92999999999
$ws.Range("K132").Value = 264965.79
$ws.Range("M132").Value = -262435.79
$ws.Range("H138").Value = 1996.683
$ws.Range("J138").Value = 3126.2307
$ws.Range("L138").Value = 9378.6921
$ws.Range("N138").Value = -19658.6921
$ws.Range("H141").Value = 1546.1666
$ws.Range("J141").Value = 1999
$ws.Range("L141").Value = 5997
$ws.Range("N141").Value = -16357

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1842
$ws.Range("I105").Value = 1862.7142
$ws.Range("J105").Value = 1769.5
$ws.Range("K105").Value = 1862.7142
$ws.Range("L105").Value = 1769.5
$ws.Range("M105").Value = -115.7141999999999
$ws.Range("N105").Value = -5263.5
$ws.Range("H134").Value = 738811.1
$ws.Range("I134").Value = 780489.75
$ws.Range("K134").Value = 2341469.25
$ws.Range("M134").Value = -2338934.25

# --- Sheet: CRP (11 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2080.8462
$ws.Range("I94").Value = 1287.6
$ws.Range("K94").Value = 1287.6
$ws.Range("M94").Value = -836.5999999999999
$ws.Range("H132").Value = 50259932
$ws.Range("I132").Value = 76926504
$ws.Range("J132").Value = 736288.5600000001
$ws.Range("K132").Value = 230779512
$ws.Range("L132").Value = 2208865.68
$ws.Range("M132").Value = -230776982
$ws.Range("N132").Value = -2213925.68

# --- Sheet: CUL (19 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 14849.615
$ws.Range("I87").Value = 10331.223
$ws.Range("K87").Value = 30993.669
$ws.Range("M87").Value = -29745.669
$ws.Range("H90").Value = 14849.615
$ws.Range("I90").Value = 10331.223
$ws.Range("K90").Value = 92981.007
$ws.Range("M90").Value = -86741.007
$ws.Range("H113").Value = 9000
$ws.Range("J113").Value = 9000
$ws.Range("L113").Value = 27000
$ws.Range("N113").Value = -31340
$ws.Range("H129").Value = 1571.8667
$ws.Range("I129").Value = 557.9
$ws.Range("J129").Value = 3599.8
$ws.Range("K129").Value = 1673.7
$ws.Range("L129").Value = 10799.4
$ws.Range("M129").Value = 3326.3
$ws.Range("N129").Value = -20799.4

# --- Sheet: GSM (34 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6098.385
$ws.Range("I70").Value = 6162.636
$ws.Range("K70").Value = 6162.636
$ws.Range("M70").Value = -5892.636
$ws.Range("H73").Value = 6098.385
$ws.Range("I73").Value = 6162.636
$ws.Range("K73").Value = 6162.636
$ws.Range("M73").Value = -5226.636
$ws.Range("H80").Value = 1708874.4
$ws.Range("I80").Value = 1708874.4
$ws.Range("K80").Value = 1708874.4
$ws.Range("M80").Value = -1707876.4
$ws.Range("H83").Value = 1708874.4
$ws.Range("I83").Value = 1708874.4
$ws.Range("K83").Value = 8544372
$ws.Range("M83").Value = -8539380
$ws.Range("H93").Value = 80249
$ws.Range("J93").Value = 80249
$ws.Range("L93").Value = 80249
$ws.Range("N93").Value = -83993
$ws.Range("H97").Value = 3159.25
$ws.Range("I97").Value = 2860.875
$ws.Range("J97").Value = 4352.75
$ws.Range("K97").Value = 2860.875
$ws.Range("L97").Value = 4352.75
$ws.Range("M97").Value = -2364.875
$ws.Range("N97").Value = -5344.75
$ws.Range("H113").Value = 2748.9167
$ws.Range("I113").Value = 1798.8
$ws.Range("J113").Value = 7499.5
$ws.Range("K113").Value = 1798.8
$ws.Range("L113").Value = 7499.5
$ws.Range("M113").Value = 371.2
$ws.Range("N113").Value = -11839.5

# --- Sheet: LTW (29 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1315
$ws.Range("I22").Value = 525.7143
$ws.Range("J22").Value = 6840
$ws.Range("K22").Value = 525.7143
$ws.Range("L22").Value = 6840
$ws.Range("M22").Value = -230.7143
$ws.Range("N22").Value = -7430
$ws.Range("H27").Value = 1315
$ws.Range("I27").Value = 525.7143
$ws.Range("J27").Value = 6840
$ws.Range("K27").Value = 525.7143
$ws.Range("L27").Value = 6840
$ws.Range("M27").Value = -418.7143
$ws.Range("N27").Value = -7054
$ws.Range("H68").Value = 4544.636
$ws.Range("I68").Value = 4873.875
$ws.Range("K68").Value = 4873.875
$ws.Range("M68").Value = -4124.875
$ws.Range("H71").Value = 4544.636
$ws.Range("I71").Value = 4873.875
$ws.Range("K71").Value = 24369.375
$ws.Range("M71").Value = -20625.375
$ws.Range("H100").Value = 19250.666
$ws.Range("I100").Value = 2625
$ws.Range("J100").Value = 52502
$ws.Range("K100").Value = 2625
$ws.Range("L100").Value = 52502
$ws.Range("M100").Value = -2084
$ws.Range("N100").Value = -53584

# --- Sheet: WVR (19 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 900000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 900000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H98").Value = 89272
$ws.Range("J98").Value = 89272
$ws.Range("L98").Value = 89272
$ws.Range("N98").Value = -95262
$ws.Range("H132").Value = 18308190
$ws.Range("I132").Value = 100615064
$ws.Range("J132").Value = 17773.777
$ws.Range("K132").Value = 301845192
$ws.Range("L132").Value = 53321.33099999999
$ws.Range("M132").Value = -301842662
$ws.Range("N132").Value = -58381.33099999999

Write-Output "Applied 169 cell updates."
